$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the day/night timeslice grouping strings ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "FaP,SaP,WaP,FaD,RaP,WaD,SaD,RaD"
$wsEv.Range("C14").Value = "RaN,RaP,FaP,SaP,SaN,WaN,WaP,FaN"

# --- Sheet "re_profiles": re-order the hydro seasonal shares (M4:O7) ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "W"
$wsRe.Range("N4").Value = 0.22555529847292916
$wsRe.Range("M5").Value = "R"
$wsRe.Range("N5").Value = 0.30301943544655252
$wsRe.Range("M6").Value = "S"
$wsRe.Range("N6").Value = 0.40439611291068944
$wsRe.Range("M7").Value = "F"
$wsRe.Range("N7").Value = 0.26702915316982878
